# Fruta / hortaliza, semanal
# Insert a new weekly record at row 150 ("Red Globe", 2022-06-02) for
# Comercializadora del Agro de Limarí - Uva, shifting the existing rows
# 150-163 down to 151-164.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 150..163 down by one row, keeping formatting (incl. the date
# number format on column D) intact.
$ws.Rows.Item(150).Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A150").Value2 = 2
$ws.Range("B150").Value2 = "Comercializadora del Agro de Limarí"
$ws.Range("C150").Value2 = "Coquimbo"
$ws.Range("D150").Value2 = 44714
$ws.Range("E150").Value2 = 4
$ws.Range("F150").Value2 = "Fruta"
$ws.Range("G150").Value2 = 100109
$ws.Range("H150").Value2 = "Uva"
$ws.Range("I150").Value2 = 100109001
$ws.Range("J150").Value2 = "Uva"
$ws.Range("K150").Value2 = "Red Globe"
$ws.Range("L150").Value2 = "Primera"
$ws.Range("M150").Value2 = 300
$ws.Range("N150").Value2 = 9000
$ws.Range("O150").Value2 = 10000
$ws.Range("P150").Value2 = 9500
$ws.Range("Q150").Value2 = "`$/bandeja 18 kilos"
$ws.Range("R150").Value2 = "Provincia de Limarí"
$ws.Range("S150").Value2 = 528
$ws.Range("T150").Value2 = 18
